# Auto-generated edit script: apply updated market/profit values per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 464.23077
$ws.Range("I12").Value = 381.1111
$ws.Range("J12").Value = 651.25
$ws.Range("K12").Value = 381.1111
$ws.Range("L12").Value = 651.25
$ws.Range("M12").Value = -211.1111
$ws.Range("N12").Value = -991.25
$ws.Range("H19").Value = 4809.9
$ws.Range("J19").Value = 4750
$ws.Range("L19").Value = 4750
$ws.Range("N19").Value = -5100
$ws.Range("H32").Value = 3562.3333
$ws.Range("I32").Value = 2000
$ws.Range("J32").Value = 3874.8
$ws.Range("K32").Value = 2000
$ws.Range("L32").Value = 3874.8
$ws.Range("M32").Value = -1674
$ws.Range("N32").Value = -4526.8
$ws.Range("H33").Value = 223.31818
$ws.Range("I33").Value = 204.22223
$ws.Range("J33").Value = 309.25
$ws.Range("K33").Value = 204.22223
$ws.Range("L33").Value = 309.25
$ws.Range("M33").Value = 24.77777
$ws.Range("N33").Value = -767.25
$ws.Range("H41").Value = 84148
$ws.Range("I41").Value = 509
$ws.Range("J41").Value = 251426
$ws.Range("K41").Value = 509
$ws.Range("L41").Value = 251426
$ws.Range("M41").Value = -69
$ws.Range("N41").Value = -252306
$ws.Range("H43").Value = 2541.6667
$ws.Range("I43").Value = 2000
$ws.Range("K43").Value = 2000
$ws.Range("M43").Value = -1931
$ws.Range("H51").Value = 4909
$ws.Range("I51").Value = 5538
$ws.Range("K51").Value = 5538
$ws.Range("M51").Value = -5054
$ws.Range("H70").Value = 16678.084
$ws.Range("I70").Value = 1630
$ws.Range("J70").Value = 27426.715
$ws.Range("K70").Value = 4890
$ws.Range("L70").Value = 82280.145
$ws.Range("M70").Value = -4620
$ws.Range("N70").Value = -82820.145
$ws.Range("H73").Value = 16678.084
$ws.Range("I73").Value = 1630
$ws.Range("J73").Value = 27426.715
$ws.Range("K73").Value = 4890
$ws.Range("L73").Value = 82280.145
$ws.Range("M73").Value = -3954
$ws.Range("N73").Value = -84152.145
$ws.Range("H74").Value = 7096.1714
$ws.Range("I74").Value = 4910
$ws.Range("K74").Value = 4910
$ws.Range("M74").Value = -3974
$ws.Range("H77").Value = 7096.1714
$ws.Range("I77").Value = 4910
$ws.Range("K77").Value = 24550
$ws.Range("M77").Value = -19870
$ws.Range("H92").Value = 1017.375
$ws.Range("I92").Value = 313.78946
$ws.Range("K92").Value = 313.78946
$ws.Range("M92").Value = 934.21054
$ws.Range("H116").Value = 161333.25
$ws.Range("J116").Value = 31000
$ws.Range("L116").Value = 31000
$ws.Range("N116").Value = -37884
$ws.Range("H132").Value = 5019.864
$ws.Range("I132").Value = 1862.8334
$ws.Range("J132").Value = 8808.299999999999
$ws.Range("K132").Value = 5588.5002
$ws.Range("L132").Value = 26424.9
$ws.Range("M132").Value = -3058.5002
$ws.Range("N132").Value = -31484.9
$ws.Range("H137").Value = 37038556
$ws.Range("I137").Value = 58824732
$ws.Range("K137").Value = 176474196
$ws.Range("M137").Value = -176471646

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 17000
$ws.Range("J24").Value = 17000
$ws.Range("L24").Value = 17000
$ws.Range("N24").Value = -17748
$ws.Range("H32").Value = 809.3823
$ws.Range("I32").Value = 777.5077
$ws.Range("K32").Value = 777.5077
$ws.Range("M32").Value = -490.5077
$ws.Range("H61").Value = 1346.7576
$ws.Range("I61").Value = 1175.2333
$ws.Range("J61").Value = 3062
$ws.Range("K61").Value = 1175.2333
$ws.Range("L61").Value = 3062
$ws.Range("M61").Value = -963.2333000000001
$ws.Range("N61").Value = -3486
$ws.Range("H100").Value = 17000
$ws.Range("J100").Value = 17000
$ws.Range("L100").Value = 17000
$ws.Range("N100").Value = -19164
$ws.Range("H132").Value = 3196.524
$ws.Range("I132").Value = 2691.1143
$ws.Range("J132").Value = 5723.5713
$ws.Range("K132").Value = 8073.342900000001
$ws.Range("L132").Value = 17170.7139
$ws.Range("M132").Value = -5543.342900000001
$ws.Range("N132").Value = -22230.7139
$ws.Range("H136").Value = 1346.7576
$ws.Range("I136").Value = 1175.2333
$ws.Range("J136").Value = 3062
$ws.Range("K136").Value = 3525.699900000001
$ws.Range("L136").Value = 9186
$ws.Range("M136").Value = -975.6999000000005
$ws.Range("N136").Value = -14286

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2742.4443
$ws.Range("I80").Value = 71
$ws.Range("K80").Value = 71
$ws.Range("M80").Value = 927
$ws.Range("H83").Value = 2742.4443
$ws.Range("I83").Value = 71
$ws.Range("K83").Value = 355
$ws.Range("M83").Value = 4637
$ws.Range("H94").Value = 1361.0667
$ws.Range("J94").Value = 1554.2
$ws.Range("L94").Value = 1554.2
$ws.Range("N94").Value = -2456.2
$ws.Range("H134").Value = 2006.4615
$ws.Range("J134").Value = 5820.5
$ws.Range("L134").Value = 17461.5
$ws.Range("N134").Value = -22531.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 952.2632
$ws.Range("J22").Value = 983.3333
$ws.Range("L22").Value = 983.3333
$ws.Range("N22").Value = -1683.3333
$ws.Range("H132").Value = 111114880
$ws.Range("I132").Value = 200001380
$ws.Range("K132").Value = 600004140
$ws.Range("M132").Value = -600001610

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 446.66666
$ws.Range("I8").Value = 446.66666
$ws.Range("K8").Value = 1339.99998
$ws.Range("M8").Value = -1200.99998
$ws.Range("H75").Value = 5467
$ws.Range("J75").Value = 8722.714
$ws.Range("L75").Value = 26168.142
$ws.Range("N75").Value = -28164.142
$ws.Range("H78").Value = 5467
$ws.Range("J78").Value = 8722.714
$ws.Range("L78").Value = 78504.42600000001
$ws.Range("N78").Value = -88488.42600000001
$ws.Range("H113").Value = 3153
$ws.Range("J113").Value = 3281.2222
$ws.Range("L113").Value = 9843.6666
$ws.Range("N113").Value = -14183.6666
$ws.Range("H137").Value = 3213.0833
$ws.Range("I137").Value = 1627.3334
$ws.Range("K137").Value = 4882.0002
$ws.Range("M137").Value = 217.9997999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 5350.4287
$ws.Range("I29").Value = 1690.8
$ws.Range("K29").Value = 1690.8
$ws.Range("M29").Value = -1400.8
$ws.Range("H70").Value = 9093
$ws.Range("I70").Value = 7389.5
$ws.Range("K70").Value = 7389.5
$ws.Range("M70").Value = -7119.5
$ws.Range("H73").Value = 9093
$ws.Range("I73").Value = 7389.5
$ws.Range("K73").Value = 7389.5
$ws.Range("M73").Value = -6453.5
$ws.Range("H113").Value = 1476.1428
$ws.Range("I113").Value = 1183.5
$ws.Range("J113").Value = 1866.3334
$ws.Range("K113").Value = 1183.5
$ws.Range("L113").Value = 1866.3334
$ws.Range("M113").Value = 986.5
$ws.Range("N113").Value = -6206.3334
$ws.Range("H132").Value = 24403686
$ws.Range("I132").Value = 35727268
$ws.Range("J132").Value = 14436.538
$ws.Range("K132").Value = 107181804
$ws.Range("L132").Value = 43309.614
$ws.Range("M132").Value = -107179274
$ws.Range("N132").Value = -48369.614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3199.8
$ws.Range("J82").Value = 2999.75
$ws.Range("L82").Value = 2999.75
$ws.Range("N82").Value = -3721.75
$ws.Range("H85").Value = 3199.8
$ws.Range("J85").Value = 2999.75
$ws.Range("L85").Value = 2999.75
$ws.Range("N85").Value = -5495.75
$ws.Range("H122").Value = 4923.5
$ws.Range("I122").Value = 4508.2
$ws.Range("K122").Value = 13524.6
$ws.Range("M122").Value = -11074.6
$ws.Range("H132").Value = 3567.5
$ws.Range("I132").Value = 2996.75
$ws.Range("J132").Value = 4328.5
$ws.Range("K132").Value = 8990.25
$ws.Range("L132").Value = 12985.5
$ws.Range("M132").Value = -6460.25
$ws.Range("N132").Value = -18045.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6700
$ws.Range("I62").Value = 9000
$ws.Range("K62").Value = 9000
$ws.Range("M62").Value = -8376
$ws.Range("H65").Value = 6700
$ws.Range("I65").Value = 9000
$ws.Range("K65").Value = 45000
$ws.Range("M65").Value = -41880
$ws.Range("H132").Value = 4655214
$ws.Range("I132").Value = 5558455.5
$ws.Range("J132").Value = 9972.429
$ws.Range("K132").Value = 16675366.5
$ws.Range("L132").Value = 29917.287
$ws.Range("M132").Value = -16672836.5
$ws.Range("N132").Value = -34977.287

Write-Host "Applied 218 cell updates across 8 sheets"